$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 684
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 684
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 684
$ws.Range("N32").Value = -1336
$ws.Range("M32").ClearContents()
$ws.Range("H113").Value = 43004.332
$ws.Range("I113").Value = 78061.84
$ws.Range("J113").Value = 1572.7273
$ws.Range("K113").Value = 78061.84
$ws.Range("L113").Value = 1572.7273
$ws.Range("M113").Value = -74807.84
$ws.Range("N113").Value = -8080.7273
$ws.Range("H129").Value = 2781.9434
$ws.Range("J129").Value = 940.21277
$ws.Range("L129").Value = 2820.63831
$ws.Range("N129").Value = -12820.63831
$ws.Range("H137").Value = 1473.5927
$ws.Range("I137").Value = 1352.2354
$ws.Range("K137").Value = 4056.7062
$ws.Range("M137").Value = -1506.7062
$ws.Range("H138").Value = 3700.2163
$ws.Range("I138").Value = 2643.8333
$ws.Range("J138").Value = 4039.7678
$ws.Range("K138").Value = 7931.499899999999
$ws.Range("L138").Value = 12119.3034
$ws.Range("M138").Value = -2791.499899999999
$ws.Range("N138").Value = -22399.3034

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 21358.857
$ws.Range("J23").Value = 9900
$ws.Range("L23").Value = 9900
$ws.Range("N23").Value = -10418
$ws.Range("H32").Value = 43547.043
$ws.Range("I32").Value = 19247.805
$ws.Range("K32").Value = 19247.805
$ws.Range("M32").Value = -18960.805
$ws.Range("H44").Value = 12996.125
$ws.Range("J44").Value = 12995.571
$ws.Range("L44").Value = 12995.571
$ws.Range("N44").Value = -13971.571
$ws.Range("H55").Value = 11957.143
$ws.Range("J55").Value = 11957.143
$ws.Range("L55").Value = 11957.143
$ws.Range("N55").Value = -12587.143
$ws.Range("H61").Value = 2659
$ws.Range("I61").Value = 2637.7
$ws.Range("J61").Value = 2730
$ws.Range("K61").Value = 2637.7
$ws.Range("L61").Value = 2730
$ws.Range("M61").Value = -2425.7
$ws.Range("N61").Value = -3154
$ws.Range("H80").Value = 26473.555
$ws.Range("J80").Value = 26473.555
$ws.Range("L80").Value = 26473.555
$ws.Range("N80").Value = -28469.555
$ws.Range("H83").Value = 26473.555
$ws.Range("J83").Value = 26473.555
$ws.Range("L83").Value = 79420.66500000001
$ws.Range("N83").Value = -89404.66500000001
$ws.Range("H122").Value = 2223.1052
$ws.Range("I122").Value = 2063.9333
$ws.Range("J122").Value = 2820
$ws.Range("K122").Value = 6191.7999
$ws.Range("L122").Value = 8460
$ws.Range("M122").Value = -3741.7999
$ws.Range("N122").Value = -13360
$ws.Range("H136").Value = 2659
$ws.Range("I136").Value = 2637.7
$ws.Range("J136").Value = 2730
$ws.Range("K136").Value = 7913.099999999999
$ws.Range("L136").Value = 8190
$ws.Range("M136").Value = -5363.099999999999
$ws.Range("N136").Value = -13290

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 558
$ws.Range("I22").Value = 563.3333
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 563.3333
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = -390.3333
$ws.Range("N22").Value = -896
$ws.Range("H82").Value = 15928.667
$ws.Range("J82").Value = 28693.334
$ws.Range("L82").Value = 28693.334
$ws.Range("N82").Value = -29459.334
$ws.Range("H85").Value = 15928.667
$ws.Range("J85").Value = 28693.334
$ws.Range("L85").Value = 28693.334
$ws.Range("N85").Value = -31345.334
$ws.Range("H120").Value = 33863.5
$ws.Range("J120").Value = 33863.5
$ws.Range("L120").Value = 33863.5
$ws.Range("N120").Value = -43539.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 908.2273
$ws.Range("I16").Value = 712
$ws.Range("J16").Value = 1251.625
$ws.Range("K16").Value = 712
$ws.Range("L16").Value = 1251.625
$ws.Range("M16").Value = -425
$ws.Range("N16").Value = -1825.625
$ws.Range("H22").Value = 463.73334
$ws.Range("I22").Value = 359.14285
$ws.Range("J22").Value = 555.25
$ws.Range("K22").Value = 359.14285
$ws.Range("L22").Value = 555.25
$ws.Range("M22").Value = -9.14285000000001
$ws.Range("N22").Value = -1255.25
$ws.Range("H31").Value = 90294.94
$ws.Range("I31").Value = 2500
$ws.Range("J31").Value = 102000.93
$ws.Range("K31").Value = 2500
$ws.Range("L31").Value = 102000.93
$ws.Range("M31").Value = -2205
$ws.Range("N31").Value = -102590.93
$ws.Range("H34").Value = 90294.94
$ws.Range("I34").Value = 2500
$ws.Range("J34").Value = 102000.93
$ws.Range("K34").Value = 2500
$ws.Range("L34").Value = 102000.93
$ws.Range("M34").Value = -2298
$ws.Range("N34").Value = -102404.93
$ws.Range("H58").Value = 1155.3405
$ws.Range("I58").Value = 986.4358999999999
$ws.Range("J58").Value = 1978.75
$ws.Range("K58").Value = 986.4358999999999
$ws.Range("L58").Value = 1978.75
$ws.Range("M58").Value = -783.4358999999999
$ws.Range("N58").Value = -2384.75
$ws.Range("H113").Value = 908.2273
$ws.Range("I113").Value = 712
$ws.Range("J113").Value = 1251.625
$ws.Range("K113").Value = 712
$ws.Range("L113").Value = 1251.625
$ws.Range("M113").Value = 1458
$ws.Range("N113").Value = -5591.625
$ws.Range("H135").Value = 46452.25
$ws.Range("J135").Value = 46452.25
$ws.Range("L135").Value = 46452.25
$ws.Range("N135").Value = -56592.25
$ws.Range("H136").Value = 1155.3405
$ws.Range("I136").Value = 986.4358999999999
$ws.Range("J136").Value = 1978.75
$ws.Range("K136").Value = 2959.3077
$ws.Range("L136").Value = 5936.25
$ws.Range("M136").Value = -409.3076999999998
$ws.Range("N136").Value = -11036.25
$ws.Range("H138").Value = 149000
$ws.Range("J138").Value = 149000
$ws.Range("L138").Value = 149000
$ws.Range("N138").Value = -159280
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 93517.336
$ws.Range("I141").Value = 92220.8
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 92220.8
$ws.Range("L141").Value = 100000
$ws.Range("M141").Value = -87040.8
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 30000
$ws.Range("N17").Value = -30338
$ws.Range("H38").Value = 195.25
$ws.Range("I38").Value = 45.5
$ws.Range("J38").Value = 245.16667
$ws.Range("K38").Value = 136.5
$ws.Range("L38").Value = 735.50001
$ws.Range("M38").Value = 210.5
$ws.Range("N38").Value = -1429.50001
$ws.Range("H113").Value = 1226
$ws.Range("I113").Value = 1704.2222
$ws.Range("J113").Value = 688
$ws.Range("K113").Value = 5112.6666
$ws.Range("L113").Value = 2064
$ws.Range("M113").Value = -2942.6666
$ws.Range("N113").Value = -6404
$ws.Range("H122").Value = 941.6923
$ws.Range("J122").Value = 1410.3334
$ws.Range("L122").Value = 12693.0006
$ws.Range("N122").Value = -17593.0006
$ws.Range("H131").Value = 714893.6
$ws.Range("I131").Value = 750
$ws.Range("J131").Value = 757529.0600000001
$ws.Range("K131").Value = 2250
$ws.Range("L131").Value = 2272587.18
$ws.Range("M131").Value = 2790
$ws.Range("N131").Value = -2282667.18
$ws.Range("H137").Value = 44738.074
$ws.Range("I137").Value = 93941.73
$ws.Range("J137").Value = 10910.5625
$ws.Range("K137").Value = 281825.19
$ws.Range("L137").Value = 32731.6875
$ws.Range("M137").Value = -276725.19
$ws.Range("N137").Value = -42931.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 75000
$ws.Range("J111").Value = 75000
$ws.Range("L111").Value = 75000
$ws.Range("N111").Value = -81134
$ws.Range("H122").Value = 2063.087
$ws.Range("I122").Value = 1758.2941
$ws.Range("J122").Value = 2926.6667
$ws.Range("K122").Value = 5274.8823
$ws.Range("L122").Value = 8780.000100000001
$ws.Range("M122").Value = -2824.8823
$ws.Range("N122").Value = -13680.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1400
$ws.Range("I22").Value = 9800
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 9800
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -9505
$ws.Range("N22").Value = -1290
$ws.Range("H27").Value = 1400
$ws.Range("I27").Value = 9800
$ws.Range("J27").Value = 700
$ws.Range("K27").Value = 9800
$ws.Range("L27").Value = 700
$ws.Range("M27").Value = -9693
$ws.Range("N27").Value = -914
$ws.Range("H40").Value = 61068.176
$ws.Range("I40").Value = 201399.8
$ws.Range("J40").Value = 2596.6667
$ws.Range("K40").Value = 201399.8
$ws.Range("L40").Value = 2596.6667
$ws.Range("M40").Value = -201263.8
$ws.Range("N40").Value = -2868.6667
$ws.Range("H46").Value = 5518
$ws.Range("J46").Value = 4562.5
$ws.Range("L46").Value = 4562.5
$ws.Range("N46").Value = -4938.5
$ws.Range("H93").Value = 4148.125
$ws.Range("I93").Value = 4169
$ws.Range("J93").Value = 4002
$ws.Range("K93").Value = 4169
$ws.Range("L93").Value = 4002
$ws.Range("M93").Value = -2921
$ws.Range("N93").Value = -6498
$ws.Range("H110").Value = 29000
$ws.Range("J110").Value = 29000
$ws.Range("L110").Value = 29000
$ws.Range("N110").Value = -37180
$ws.Range("H132").Value = 4173.875
$ws.Range("J132").Value = 2766.2727
$ws.Range("L132").Value = 8298.8181
$ws.Range("N132").Value = -13358.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 22607.22
$ws.Range("I132").Value = 2869.3823
$ws.Range("J132").Value = 64550.125
$ws.Range("K132").Value = 8608.1469
$ws.Range("L132").Value = 193650.375
$ws.Range("M132").Value = -6078.1469
$ws.Range("N132").Value = -198710.375
